$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 31250900
$ws.Range("I33").Value = 636.2143
$ws.Range("J33").Value = 250002750
$ws.Range("K33").Value = 636.2143
$ws.Range("L33").Value = 250002750
$ws.Range("M33").Value = -407.2143
$ws.Range("N33").Value = -250003208

$ws.Range("H135").Value = 877.3333
$ws.Range("I135").Value = 877.3333
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 7895.9997
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = ""
$ws.Range("N135").Value = -5360.9997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H114").Value = 30398
$ws.Range("J114").Value = 30398
$ws.Range("L114").Value = 30398
$ws.Range("N114").Value = -39076

$ws.Range("H132").Value = 2174.4644
$ws.Range("I132").Value = 1741.5714
$ws.Range("J132").Value = 3473.1428
$ws.Range("K132").Value = 5224.7142
$ws.Range("L132").Value = 10419.4284
$ws.Range("M132").Value = -2694.7142
$ws.Range("N132").Value = -15479.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 26250
$ws.Range("J21").Value = 26250
$ws.Range("L21").Value = 26250
$ws.Range("N21").Value = -26722

$ws.Range("H37").Value = 8818.727999999999
$ws.Range("I37").Value = 3143.7144
$ws.Range("J37").Value = 18750
$ws.Range("K37").Value = 3143.7144
$ws.Range("L37").Value = 18750
$ws.Range("M37").Value = -3006.7144
$ws.Range("N37").Value = -19024

$ws.Range("H115").Value = 20000
$ws.Range("J115").Value = 20000
$ws.Range("L115").Value = 20000
$ws.Range("N115").Value = -23134

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11112563
$ws.Range("I31").Value = 1000.0476
$ws.Range("J31").Value = 20835180
$ws.Range("K31").Value = 1000.0476
$ws.Range("L31").Value = 20835180
$ws.Range("M31").Value = -705.0476
$ws.Range("N31").Value = -20835770

$ws.Range("H34").Value = 11112563
$ws.Range("I34").Value = 1000.0476
$ws.Range("J34").Value = 20835180
$ws.Range("K34").Value = 1000.0476
$ws.Range("L34").Value = 20835180
$ws.Range("M34").Value = -798.0476
$ws.Range("N34").Value = -20835584

$ws.Range("H58").Value = 5535
$ws.Range("J58").Value = 600
$ws.Range("L58").Value = 600
$ws.Range("N58").Value = -1006

$ws.Range("H86").Value = 13400.25
$ws.Range("I86").Value = 4201.4
$ws.Range("J86").Value = 28731.666
$ws.Range("K86").Value = 4201.4
$ws.Range("L86").Value = 28731.666
$ws.Range("M86").Value = -3078.4
$ws.Range("N86").Value = -30977.666

$ws.Range("H89").Value = 13400.25
$ws.Range("I89").Value = 4201.4
$ws.Range("J89").Value = 28731.666
$ws.Range("K89").Value = 21007
$ws.Range("L89").Value = 143658.33
$ws.Range("M89").Value = -15391
$ws.Range("N89").Value = -154890.33

$ws.Range("H122").Value = 3335763.2
$ws.Range("I122").Value = 2780
$ws.Range("J122").Value = 6668746.5
$ws.Range("K122").Value = 8340
$ws.Range("L122").Value = 20006239.5
$ws.Range("M122").Value = -5890
$ws.Range("N122").Value = -20011139.5

$ws.Range("H136").Value = 5535
$ws.Range("J136").Value = 600
$ws.Range("L136").Value = 1800
$ws.Range("N136").Value = -6900

$ws.Range("H138").Value = 39125.168
$ws.Range("J138").Value = 39125.168
$ws.Range("L138").Value = 39125.168
$ws.Range("N138").Value = -49405.168

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 3500.5
$ws.Range("J25").Value = 5001
$ws.Range("L25").Value = 15003
$ws.Range("N25").Value = -15341

$ws.Range("H30").Value = 3500.5
$ws.Range("J30").Value = 5001
$ws.Range("L30").Value = 15003
$ws.Range("N30").Value = -15207

$ws.Range("H35").Value = 1425.1428
$ws.Range("J35").Value = 1425.1428
$ws.Range("L35").Value = 4275.428400000001
$ws.Range("N35").Value = -4851.428400000001

$ws.Range("H68").Value = 1729.098
$ws.Range("I68").Value = 837.5599999999999
$ws.Range("J68").Value = 2586.3462
$ws.Range("K68").Value = 2512.68
$ws.Range("L68").Value = 7759.0386
$ws.Range("M68").Value = -1701.68
$ws.Range("N68").Value = -9381.0386

$ws.Range("H71").Value = 1729.098
$ws.Range("I71").Value = 837.5599999999999
$ws.Range("J71").Value = 2586.3462
$ws.Range("K71").Value = 7538.039999999999
$ws.Range("L71").Value = 23277.1158
$ws.Range("M71").Value = -3482.039999999999
$ws.Range("N71").Value = -31389.1158

$ws.Range("H107").Value = 586390.2
$ws.Range("I107").Value = 584.3200000000001
$ws.Range("J107").Value = 885270.75
$ws.Range("K107").Value = 1752.96
$ws.Range("L107").Value = 2655812.25
$ws.Range("M107").Value = 167.04
$ws.Range("N107").Value = -2659652.25

$ws.Range("H140").Value = 1815.2941
$ws.Range("I140").Value = 705
$ws.Range("J140").Value = 4480
$ws.Range("K140").Value = 2115
$ws.Range("L140").Value = 13440
$ws.Range("M140").Value = 3065
$ws.Range("N140").Value = -23800

$ws.Range("H141").Value = 2900
$ws.Range("I141").Value = 2900
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 8700
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = ""
$ws.Range("N141").Value = -3520

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1044.0741
$ws.Range("I97").Value = 914.8889
$ws.Range("J97").Value = 1302.4445
$ws.Range("K97").Value = 914.8889
$ws.Range("L97").Value = 1302.4445
$ws.Range("M97").Value = -418.8889
$ws.Range("N97").Value = -2294.4445

$ws.Range("H99").Value = 10761.833
$ws.Range("I99").Value = 7114.2
$ws.Range("J99").Value = 29000
$ws.Range("K99").Value = 7114.2
$ws.Range("L99").Value = 29000
$ws.Range("M99").Value = -4868.2
$ws.Range("N99").Value = -33492

$ws.Range("H103").Value = 24000
$ws.Range("J103").Value = 24000
$ws.Range("L103").Value = 24000
$ws.Range("N103").Value = -26344

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = ""
$ws.Range("N111").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 934.88
$ws.Range("I46").Value = 1280.2
$ws.Range("J46").Value = 848.55
$ws.Range("K46").Value = 1280.2
$ws.Range("L46").Value = 848.55
$ws.Range("M46").Value = -1092.2
$ws.Range("N46").Value = -1224.55

$ws.Range("H55").Value = 455.8095
$ws.Range("I55").Value = 204
$ws.Range("J55").Value = 1966.6666
$ws.Range("K55").Value = 204
$ws.Range("L55").Value = 1966.6666
$ws.Range("M55").Value = -31
$ws.Range("N55").Value = -2312.6666

$ws.Range("H61").Value = 2084.8948
$ws.Range("I61").Value = 1613.5
$ws.Range("J61").Value = 2427.7273
$ws.Range("K61").Value = 1613.5
$ws.Range("L61").Value = 2427.7273
$ws.Range("M61").Value = -1411.5
$ws.Range("N61").Value = -2831.7273

$ws.Range("H94").Value = 15404.143
$ws.Range("J94").Value = 15404.143
$ws.Range("L94").Value = 15404.143
$ws.Range("N94").Value = -16756.143

$ws.Range("H95").Value = 25114.666
$ws.Range("I95").Value = 23000
$ws.Range("J95").Value = 26172
$ws.Range("K95").Value = 23000
$ws.Range("L95").Value = 26172
$ws.Range("M95").Value = -20254
$ws.Range("N95").Value = -31664

$ws.Range("H96").Value = 26333.334
$ws.Range("I96").Value = 30000
$ws.Range("J96").Value = 24500
$ws.Range("K96").Value = 30000
$ws.Range("L96").Value = 24500
$ws.Range("M96").Value = -27254
$ws.Range("N96").Value = -29992

$ws.Range("H97").Value = 25000
$ws.Range("J97").Value = 25000
$ws.Range("L97").Value = 25000
$ws.Range("N97").Value = -26982

$ws.Range("H99").Value = 28000
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").Value = ""

$ws.Range("H109").Value = 28150
$ws.Range("J109").Value = 28150
$ws.Range("L109").Value = 28150
$ws.Range("N109").Value = -30924

$ws.Range("H113").Value = 2084.8948
$ws.Range("I113").Value = 1613.5
$ws.Range("J113").Value = 2427.7273
$ws.Range("K113").Value = 1613.5
$ws.Range("L113").Value = 2427.7273
$ws.Range("M113").Value = 556.5
$ws.Range("N113").Value = -6767.7273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 20000
$ws.Range("J40").Value = 20000
$ws.Range("L40").Value = 20000
$ws.Range("N40").Value = -20298

$ws.Range("H94").Value = 29000
$ws.Range("J94").Value = 29000
$ws.Range("L94").Value = 29000
$ws.Range("N94").Value = -30802

$ws.Range("H113").Value = 465.91306
$ws.Range("I113").Value = 362.69232
$ws.Range("J113").Value = 600.1
$ws.Range("K113").Value = 1088.07696
$ws.Range("L113").Value = 1800.3
$ws.Range("M113").Value = 1081.92304
$ws.Range("N113").Value = -6140.3
